$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: "6. " | "Zasilacz 12V" | 1 ---
$ws.Range("A7").Value = "'6. "
$ws.Range("B7").Value = "Zasilacz 12V"
$ws.Range("C7").Value = 1

# --- New row 8: "7." | "Przewody" ---
$ws.Range("A8").Value = "'7."
$ws.Range("B8").Value = "Przewody"

# --- New notes appended to row 5 ---
$ws.Range("F5").Value = "tranzystorek sterowany z raspberry"
$ws.Range("J5").Value = "przekaźnik"
$ws.Range("E5").Value = "> zamiast"

# --- New row 9: "8. " | "Przekaźnik" | 2 | "Przekaznik" (hyperlink) ---
$ws.Range("A9").Value = "'8. "
$ws.Range("D9").Value = "Przekaznik"
$ws.Hyperlinks.Add($ws.Range("D9"), "https://botland.com.pl/przekazniki/7469-przekaznik-1-stykowy-12v-10a-z-gniazdem.html")
$ws.Range("D9").Style = "Hiperłącze"
$ws.Range("B9").Value = "Przekaźnik"
$ws.Range("C9").Value = 2

# --- More notes on row 5 ---
$ws.Range("K5").Value = " Możliwe że przetwornica jest nie potrzebna"

# --- Final selection matches authored state ---
$ws.Range("M8").Select()
